# Scheduled runner refresh: update cached market-board figures
# (currentAveragePrice* / LevePrice* / LeveProfit* columns) per sheet.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H88").Value = 1590.7693
$ws.Range("I88").Value = 1550.6
$ws.Range("K88").Value = 1550.6
$ws.Range("M88").Value = -1144.6

$ws.Range("H91").Value = 1590.7693
$ws.Range("I91").Value = 1550.6
$ws.Range("K91").Value = 1550.6
$ws.Range("M91").Value = -146.5999999999999

$ws.Range("H116").Value = 11016.65
$ws.Range("I116").Value = 13634.846
$ws.Range("J116").Value = 6154.2856
$ws.Range("K116").Value = 13634.846
$ws.Range("L116").Value = 6154.2856
$ws.Range("M116").Value = -10192.846
$ws.Range("N116").Value = -13038.2856

$ws.Range("H127").Value = 24525.691
$ws.Range("I127").Value = 28439.818
$ws.Range("J127").Value = 2998
$ws.Range("K127").Value = 85319.454
$ws.Range("L127").Value = 8994
$ws.Range("M127").Value = -80359.454
$ws.Range("N127").Value = -18914

$ws.Range("H132").Value = 1283.4584
$ws.Range("I132").Value = 1315.9025
$ws.Range("K132").Value = 3947.7075
$ws.Range("M132").Value = -1417.7075

$ws.Range("H135").Value = 1689.7059
$ws.Range("I135").Value = 1664.3572
$ws.Range("J135").Value = 1808
$ws.Range("K135").Value = 14979.2148
$ws.Range("L135").Value = 16272
$ws.Range("M135").Value = -12444.2148
$ws.Range("N135").Value = -21342

$ws.Range("H137").Value = 2663664.2
$ws.Range("I137").Value = 3613100.2
$ws.Range("J137").Value = 5243.8
$ws.Range("K137").Value = 10839300.6
$ws.Range("L137").Value = 15731.4
$ws.Range("M137").Value = -10836750.6
$ws.Range("N137").Value = -20831.4

$ws.Range("H138").Value = 4082.13
$ws.Range("I138").Value = 2099.9565
$ws.Range("J138").Value = 4674.208
$ws.Range("K138").Value = 6299.869499999999
$ws.Range("L138").Value = 14022.624
$ws.Range("M138").Value = -1159.869499999999
$ws.Range("N138").Value = -24302.624

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1959
$ws.Range("J45").Value = 2500
$ws.Range("L45").Value = 2500
$ws.Range("N45").Value = -3254

$ws.Range("H88").Value = 1113.4375
$ws.Range("I88").Value = 1306.1
$ws.Range("J88").Value = 1025.8636
$ws.Range("K88").Value = 1306.1
$ws.Range("L88").Value = 1025.8636
$ws.Range("M88").Value = -900.0999999999999
$ws.Range("N88").Value = -1837.8636

$ws.Range("H91").Value = 1113.4375
$ws.Range("I91").Value = 1306.1
$ws.Range("J91").Value = 1025.8636
$ws.Range("K91").Value = 1306.1
$ws.Range("L91").Value = 1025.8636
$ws.Range("M91").Value = 97.90000000000009
$ws.Range("N91").Value = -3833.8636

$ws.Range("H97").Value = 1684.4445
$ws.Range("J97").Value = 2808.8572
$ws.Range("L97").Value = 2808.8572
$ws.Range("N97").Value = -3800.8572

$ws.Range("H104").Value = 32000
$ws.Range("J104").Value = 32000
$ws.Range("L104").Value = 32000
$ws.Range("N104").Value = -38988

$ws.Range("H122").Value = 12217.733
$ws.Range("I122").Value = 7772.636
$ws.Range("K122").Value = 23317.908
$ws.Range("M122").Value = -20867.908

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2195.7273
$ws.Range("I105").Value = 2265.4
$ws.Range("J105").Value = 1499
$ws.Range("K105").Value = 2265.4
$ws.Range("L105").Value = 1499
$ws.Range("M105").Value = -518.4000000000001
$ws.Range("N105").Value = -4993

$ws.Range("H134").Value = 2330.9143
$ws.Range("I134").Value = 2021.3214
$ws.Range("J134").Value = 3569.2856
$ws.Range("K134").Value = 6063.9642
$ws.Range("L134").Value = 10707.8568
$ws.Range("M134").Value = -3528.9642
$ws.Range("N134").Value = -15777.8568

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 569.7
$ws.Range("I22").Value = 388.66666
$ws.Range("K22").Value = 388.66666
$ws.Range("M22").Value = -38.66665999999998

$ws.Range("H31").Value = 270558
$ws.Range("I31").Value = 334829.8
$ws.Range("J31").Value = 77742.60000000001
$ws.Range("K31").Value = 334829.8
$ws.Range("L31").Value = 77742.60000000001
$ws.Range("M31").Value = -334534.8
$ws.Range("N31").Value = -78332.60000000001

$ws.Range("H34").Value = 270558
$ws.Range("I34").Value = 334829.8
$ws.Range("J34").Value = 77742.60000000001
$ws.Range("K34").Value = 334829.8
$ws.Range("L34").Value = 77742.60000000001
$ws.Range("M34").Value = -334627.8
$ws.Range("N34").Value = -78146.60000000001

$ws.Range("H58").Value = 3339.4285
$ws.Range("I58").Value = 3579.5
$ws.Range("K58").Value = 3579.5
$ws.Range("M58").Value = -3376.5

$ws.Range("H105").Value = 4459.85
$ws.Range("I105").Value = 1871.2142
$ws.Range("J105").Value = 5853.731
$ws.Range("K105").Value = 1871.2142
$ws.Range("L105").Value = 5853.731
$ws.Range("M105").Value = -124.2141999999999
$ws.Range("N105").Value = -9347.731

$ws.Range("H107").Value = 4717.2
$ws.Range("I107").Value = 668.61536
$ws.Range("J107").Value = 7109.5454
$ws.Range("K107").Value = 668.61536
$ws.Range("L107").Value = 7109.5454
$ws.Range("M107").Value = 1251.38464
$ws.Range("N107").Value = -10949.5454

$ws.Range("H134").Value = 8002.6665
$ws.Range("I134").Value = 9062.467000000001
$ws.Range("J134").Value = 2703.6667
$ws.Range("K134").Value = 27187.401
$ws.Range("L134").Value = 8111.000100000001
$ws.Range("M134").Value = -24652.401
$ws.Range("N134").Value = -13181.0001

$ws.Range("H136").Value = 3339.4285
$ws.Range("I136").Value = 3579.5
$ws.Range("K136").Value = 10738.5
$ws.Range("M136").Value = -8188.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H104").Value = 32500
$ws.Range("J104").Value = 32500
$ws.Range("L104").Value = 32500
$ws.Range("N104").Value = -39488

$ws.Range("H107").Value = 806.0625
$ws.Range("I107").Value = 425.25
$ws.Range("J107").Value = 1948.5
$ws.Range("K107").Value = 425.25
$ws.Range("L107").Value = 1948.5
$ws.Range("M107").Value = 1494.75
$ws.Range("N107").Value = -5788.5

$ws.Range("H126").Value = 11726.833
$ws.Range("I126").Value = 40000
$ws.Range("K126").Value = 120000
$ws.Range("M126").Value = -117530

$ws.Range("H132").Value = 112240
$ws.Range("I132").Value = 181668.67
$ws.Range("J132").Value = 42811.332
$ws.Range("K132").Value = 545006.01
$ws.Range("L132").Value = 128433.996
$ws.Range("M132").Value = -542476.01
$ws.Range("N132").Value = -133493.996

$ws.Range("H136").Value = 19051.84
$ws.Range("J136").Value = 19051.84
$ws.Range("L136").Value = 57155.52
$ws.Range("N136").Value = -62255.52

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H64").Value = 34284
$ws.Range("J64").Value = 34284
$ws.Range("L64").Value = 34284
$ws.Range("N64").Value = -34734

$ws.Range("H67").Value = 34284
$ws.Range("J67").Value = 34284
$ws.Range("L67").Value = 34284
$ws.Range("N67").Value = -35844

$ws.Range("H132").Value = 6481.4385
$ws.Range("I132").Value = 3215.9644
$ws.Range("K132").Value = 9647.893199999999
$ws.Range("M132").Value = -7117.893199999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H31").Value = 13666.333
$ws.Range("J31").Value = 14999.5
$ws.Range("L31").Value = 14999.5
$ws.Range("N31").Value = -15695.5

$ws.Range("H132").Value = 1382.8695
$ws.Range("I132").Value = 1268.4546
$ws.Range("J132").Value = 3900
$ws.Range("K132").Value = 3805.3638
$ws.Range("L132").Value = 11700
$ws.Range("M132").Value = -1275.3638
$ws.Range("N132").Value = -16760
